$wb = $excel.ActiveWorkbook
Write-Host ($wb | Get-Member)
